$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 test data ----------------------------------------------------
# Title / FirstName / Surname (existing shared strings get re-pointed)
$ws.Range("B2").Value = "Lieutenant Colonel"
$ws.Range("C2").Value = "Gowtham"
$ws.Range("D2").Value = "S"

# Provision before Month so the new shared strings land in the same
# order as the target workbook (Provisional=22, May=23).
$ws.Range("I2").Value = "Provisional"
$ws.Range("G2").Value = "May"

# Year / Date are plain numbers. The cells are formatted as Text ("@"),
# which would normally coerce a typed value into a string, so flip the
# number format to a non-text one for the write and back again.
$ws.Range("F2").NumberFormat = "0"
$ws.Range("F2").Value = 1998
$ws.Range("F2").NumberFormat = "@"

$ws.Range("H2").NumberFormat = "0"
$ws.Range("H2").Value = 30
$ws.Range("H2").NumberFormat = "@"

# --- Column B width -------------------------------------------------------
# Excel stores column widths in pixels internally and re-derives the
# character-based ColumnWidth from that, so asking for the nominal target
# (23) overshoots after the round-trip; 22.14 is what lands on an OOXML
# <col width="23"/> after that pixel-snapping.
$ws.Columns("B").ColumnWidth = 22.14

# --- Selection --------------------------------------------------------------
$ws.Range("F6").Select()

# --- Data validations ---------------------------------------------------------
# Remove the old date-range validation that lived on H2
$ws.Range("H2").Validation.Delete()

# Title list validation on B2
$ws.Range("B2").Validation.Add(3, 1, 1, """Mr,Mrs,Miss,Ms,Doctor,Captain,Duchess,Duke,Father,General,Lady,Lord,Lieutenant,Lieutenant Colonel,Major,Master,Professor,Reverend,Sir,Squire,Squadron Leader""")

# Provision list validation on I2
$ws.Range("I2").Validation.Add(3, 1, 1, """Full,Provisional""")
